# daily auto push: 2026-01-08 02:26 UTC
#
# The log sheet gained one more reading for 2026/01/08 (time=10) that was
# missing between the existing 06:00 entry and the next block (2026/12/29).
# Insert a new row 591 and push every following row (old 591..632) down by
# one (new 592..633), matching the diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 591:632 down to 592:633, opening up a blank row 591.
$ws.Rows.Item(591).Insert()

# Column A stores the date as literal text (e.g. "2026/01/08"), not a real
# date cell -- the sheet has no date number format anywhere. Prefix with an
# apostrophe so Excel stores it as text instead of auto-converting it to a
# date serial, then clear the resulting cell formatting so no stray
# "quote prefix" style sticks around (matching the unstyled cells around it).
$ws.Cells.Item(591, 1).Value = "'2026/01/08"
$ws.Cells.Item(591, 1).ClearFormats()
$ws.Cells.Item(591, 2).Value = "木"
$ws.Cells.Item(591, 3).Value = 10
$ws.Cells.Item(591, 4).Value = 201
